# Insert a new "Force" row above the existing "SqlName" row (row 8) in the
# Property1 sheet, pushing it and all rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# 1. Insert a blank row at row 8 (existing rows 8-10 shift down to 9-11)
$ws.Rows("8:8").Insert()

# 2. Copy the formatting of the row directly below (the row that used to be
#    row 8, now row 9) into the newly inserted row 8, so the new row matches
#    the look of its neighbours exactly.
$ws.Range("A9:I9").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Fill in the new row's content: label "Force" in column A and FALSE
#    boolean flags across columns B:I (same pattern as the other flag rows).
$ws.Range("A8").Value2 = "Force"
$ws.Range("B8:I8").Value2 = $false

# 4. Re-establish the frozen pane so the split moves from after row 9 to
#    after row 10 (i.e. below the newly inserted row), matching the shift of
#    all the data caused by the insert.
$win.FreezePanes = $false
[void]$ws.Range("A11").Select()
$win.FreezePanes = $true

# 5. Restore the selected/active cell.
[void]$ws.Range("A9").Select()
